$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row right below the header row to hold the newest meeting (第614回).
$ws.Rows.Item(2).Insert()

$topic = @"
１入院・外来医療等の調査・評価分科会からの報告について
２医薬品の新規薬価収載について
３医療機器及び臨床検査の保険適用について
４ＤＰＣにおける高額な新規の医薬品等への対応について
５保険医が投与することができる注射薬について
６令和８年度診療報酬改定におけるＤＰＣ制度への参加又はＤＰＣ制度からの退出に係る届出の受付期間について
７費用対効果評価の結果を踏まえた薬価の見直しについて
８高額医薬品（認知症薬）に対する対応について

"@

$material = @"
資料

"@

$ws.Range("A2").Value = "第614回"
$ws.Range("B2").Value = "2025年8月6日（令和7年8月6日）"
$ws.Range("C2").Value = $topic
$ws.Range("D2").Value = "－"
$ws.Range("E2").Value = $material
$ws.Range("F2").Value = "－"
